$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 05:46:45"
$wsZhCn.Range("H2").Value = "2016-03-22 05:47:23"

# "de-de" sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 05:46:52"
$wsDeDe.Range("H2").Value = "2016-03-22 05:47:36"
